$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("2023-12-07 12:02:09", 0.001),
    @("2023-12-07 12:02:36", 0.0012),
    @("2023-12-07 12:03:25", 0.0034),
    @("2023-12-07 12:03:38", 0.001),
    @("2023-12-07 12:04:08", 0.001)
)

$startRow = 64
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
